$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update the "time_taken" timestamps on the data sheet (re-run at a later time)
$ws.Range("F2").Value = "2021-10-05 14:19:26.986046"
$ws.Range("F3").Value = "2021-10-05 14:19:26.986054"
$ws.Range("F4").Value = "2021-10-05 14:19:26.986058"
$ws.Range("F5").Value = "2021-10-05 14:19:26.986060"
$ws.Range("F6").Value = "2021-10-05 14:19:26.986063"
$ws.Range("F7").Value = "2021-10-05 14:19:26.986066"
$ws.Range("F8").Value = "2021-10-05 14:19:26.986068"
$ws.Range("F9").Value = "2021-10-05 14:19:26.986071"
$ws.Range("F10").Value = "2021-10-05 14:19:26.986074"
$ws.Range("F11").Value = "2021-10-05 14:19:26.986077"

# Add the new "metadata" sheet right after "data"
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Copy the header style (bold + border, same as data!B1) onto the new header row
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the data-row style (bordered, same as data!A2) onto the new A2 cell
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Catecholaminergic polymorphic VT"
$meta.Range("C2").Value = 214

# data_version must stay text ("2.19"), not be coerced to a number - force
# text storage, then strip the leftover number-format so the cell keeps the
# default (unstyled) look by re-pasting an unformatted neighbour's format.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.19"
$meta.Range("C2").Copy()
$meta.Range("D2").PasteSpecial(-4122)

$meta.Range("E2").Value = "2021-09-28T09:49:49.887178Z"
$meta.Range("F2").Value = "2021-10-05 14:19:26.982651"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/214/?format=json"

# Keep "data" as the active sheet/selection, as in the original workbook
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
